# Training Record.xlsx edit:
# - Add new worksheet "After csf 143 games" (143-game training data),
#   built as a duplicate of "After 98 Games" so it inherits the same
#   per-column number fills/styles, then filled in with its own counts.
# - The new sheet becomes the selected/active tab.
# - "After 98 Games" loses its tab-selected flag; its recorded selection
#   becomes a whole-sheet selection (A1:XFD1048576).

$wb = $excel.ActiveWorkbook

# "After 98 Games" is the template: identical column layout/styles to the new sheet.
$wsTemplate = $wb.Worksheets.Item("After 98 Games")
$wsTemplate.Cells.Select() | Out-Null

# Duplicate it (this preserves column styles/fills/widths) and rename the copy.
$wsTemplate.Copy($null, $wsTemplate)
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "After csf 143 games"

# Overwrite with the "After csf 143 games" box-fill counts.
$ws3.Range("A2").Value = 2
$ws3.Range("B2").Value = 5
$ws3.Range("A3").Value = 3
$ws3.Range("B3").Value = 24
$ws3.Range("C3").Value = 1
$ws3.Range("A4").Value = 4
$ws3.Range("B4").Value = 31
$ws3.Range("C4").Value = 1
$ws3.Range("D4").Value = 2
$ws3.Range("A5").Value = 5
$ws3.Range("E5").Value = 10
$ws3.Range("A6").Value = 6
$ws3.Range("B6").Value = 2
$ws3.Range("E6").Value = 2
$ws3.Range("A7").Value = 7
$ws3.Range("C7").Value = 17
$ws3.Range("E7").Value = 1
$ws3.Range("A8").Value = 8
$ws3.Range("C8").Value = 8
$ws3.Range("D8").Value = 1
$ws3.Range("E8").Value = 1
$ws3.Range("A9").Value = 9
$ws3.Range("F9").Value = 50
$ws3.Range("A10").Value = 10
$ws3.Range("E10").Value = 1
$ws3.Range("F10").Value = 2
$ws3.Range("A11").Value = 11
$ws3.Range("B11").Value = 2
$ws3.Range("C11").Value = 1
$ws3.Range("D11").Value = 8
$ws3.Range("E11").Value = 2
$ws3.Range("F11").Value = 1
$ws3.Range("A12").Value = 12
$ws3.Range("B12").Value = 1
$ws3.Range("C12").Value = 1
$ws3.Range("E12").Value = 1
$ws3.Range("F12").Value = 3
$ws3.Range("G12").Value = 2
$ws3.Range("A13").Value = 13
$ws3.Range("B13").Value = 1
$ws3.Range("C13").Value = 1
$ws3.Range("D13").Value = 2
$ws3.Range("E13").Value = 3
$ws3.Range("F13").Value = 2
$ws3.Range("G13").Value = 1
$ws3.Range("A14").Value = 14
$ws3.Range("B14").Value = 3
$ws3.Range("C14").Value = 1
$ws3.Range("D14").Value = 2
$ws3.Range("E14").Value = 2
$ws3.Range("F14").Value = 2
$ws3.Range("G14").Value = 3
$ws3.Range("H14").Value = 6
$ws3.Range("A15").Value = 15
$ws3.Range("E15").Value = 20
$ws3.Range("A16").Value = 16
$ws3.Range("B16").Value = 1
$ws3.Range("I16").Value = 47
$ws3.Range("A17").Value = 17
$ws3.Range("I17").Value = 3
$ws3.Range("A18").Value = 18
$ws3.Range("C18").Value = 5
$ws3.Range("D18").Value = 19
$ws3.Range("I18").Value = 8
$ws3.Range("A19").Value = 19
$ws3.Range("I19").Value = 4
$ws3.Range("A20").Value = 20
$ws3.Range("B20").Value = 2
$ws3.Range("C20").Value = 2
$ws3.Range("E20").Value = 2
$ws3.Range("F20").Value = 1
$ws3.Range("I20").Value = 9
$ws3.Range("A21").Value = 21
$ws3.Range("B21").Value = 2
$ws3.Range("C21").Value = 5
$ws3.Range("D21").Value = 3
$ws3.Range("F21").Value = 6
$ws3.Range("I21").Value = 4
$ws3.Range("A22").Value = 22
$ws3.Range("C22").Value = 2
$ws3.Range("E22").Value = 2
$ws3.Range("F22").Value = 2
$ws3.Range("G22").Value = 1
$ws3.Range("I22").Value = 2
$ws3.Range("A23").Value = 23
$ws3.Range("B23").Value = 2
$ws3.Range("C23").Value = 4
$ws3.Range("D23").Value = 2
$ws3.Range("E23").Value = 2
$ws3.Range("F23").Value = 2
$ws3.Range("G23").Value = 2
$ws3.Range("I23").Value = 2
$ws3.Range("A24").Value = 24
$ws3.Range("B24").Value = 1
$ws3.Range("C24").Value = 4
$ws3.Range("D24").Value = 6
$ws3.Range("E24").Value = 4
$ws3.Range("F24").Value = 4
$ws3.Range("G24").Value = 3
$ws3.Range("H24").Value = 2
$ws3.Range("I24").Value = 4
$ws3.Range("A25").Value = 25
$ws3.Range("B25").Value = 3
$ws3.Range("J25").Value = 35
$ws3.Range("A26").Value = 26
$ws3.Range("C26").Value = 1
$ws3.Range("E26").Value = 1
$ws3.Range("F26").Value = 10
$ws3.Range("I26").Value = 3
$ws3.Range("A27").Value = 27
$ws3.Range("C27").Value = 11
$ws3.Range("D27").Value = 2
$ws3.Range("E27").Value = 1
$ws3.Range("F27").Value = 18
$ws3.Range("I27").Value = 1
$ws3.Range("J27").Value = 6
$ws3.Range("A28").Value = 28
$ws3.Range("B28").Value = 1
$ws3.Range("C28").Value = 1
$ws3.Range("E28").Value = 2
$ws3.Range("F28").Value = 4
$ws3.Range("G28").Value = 2
$ws3.Range("I28").Value = 4
$ws3.Range("J28").Value = 2
$ws3.Range("A29").Value = 29
$ws3.Range("B29").Value = 2
$ws3.Range("C29").Value = 2
$ws3.Range("D29").Value = 2
$ws3.Range("E29").Value = 2
$ws3.Range("F29").Value = 2
$ws3.Range("G29").Value = 2
$ws3.Range("I29").Value = 4
$ws3.Range("J29").Value = 2
$ws3.Range("A30").Value = 30
$ws3.Range("B30").Value = 1
$ws3.Range("C30").Value = 4
$ws3.Range("D30").Value = 3
$ws3.Range("F30").Value = 2
$ws3.Range("G30").Value = 3
$ws3.Range("H30").Value = 5
$ws3.Range("I30").Value = 8
$ws3.Range("J30").Value = 2
$ws3.Range("A31").Value = 31
$ws3.Range("C31").Value = 1
$ws3.Range("E31").Value = 1
$ws3.Range("F31").Value = 4
$ws3.Range("G31").Value = 11
$ws3.Range("I31").Value = 7
$ws3.Range("J31").Value = 4
$ws3.Range("K31").Value = 3
$ws3.Range("A32").Value = 32
$ws3.Range("B32").Value = 2
$ws3.Range("C32").Value = 2
$ws3.Range("D32").Value = 4
$ws3.Range("E32").Value = 1
$ws3.Range("F32").Value = 1
$ws3.Range("G32").Value = 5
$ws3.Range("I32").Value = 2
$ws3.Range("J32").Value = 2
$ws3.Range("K32").Value = 2
$ws3.Range("A33").Value = 33
$ws3.Range("B33").Value = 1
$ws3.Range("C33").Value = 1
$ws3.Range("D33").Value = 1
$ws3.Range("E33").Value = 1
$ws3.Range("F33").Value = 7
$ws3.Range("G33").Value = 18
$ws3.Range("H33").Value = 4
$ws3.Range("J33").Value = 3
$ws3.Range("K33").Value = 4
$ws3.Range("A34").Value = 34
$ws3.Range("C34").Value = 31
$ws3.Range("D34").Value = 11
$ws3.Range("E34").Value = 1
$ws3.Range("H34").Value = 19
$ws3.Range("I34").Value = 3
$ws3.Range("J34").Value = 2
$ws3.Range("K34").Value = 13
$ws3.Range("L34").Value = 25

# Cells that are populated on the template but blank on this sheet.
$ws3.Range("E9").Clear()
$ws3.Range("B10").Clear()
$ws3.Range("C10").Clear()
$ws3.Range("I15").Clear()
$ws3.Range("E16").Clear()
$ws3.Range("B17").Clear()
$ws3.Range("C17").Clear()
$ws3.Range("E17").Clear()
$ws3.Range("E19").Clear()
$ws3.Range("E21").Clear()
$ws3.Range("B22").Clear()
$ws3.Range("I25").Clear()
$ws3.Range("B26").Clear()
$ws3.Range("J26").Clear()
$ws3.Range("B27").Clear()
$ws3.Range("E30").Clear()
$ws3.Range("B31").Clear()
$ws3.Range("I33").Clear()
$ws3.Range("F34").Clear()
$ws3.Range("G34").Clear()

# Make the new sheet the active / selected tab, matching the recorded selection.
$ws3.Select() | Out-Null
$ws3.Range("B31").Select() | Out-Null

